$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("M2").Value = 45.1830845
$ws.Range("N2").Value = 90.366169
$ws.Range("O2").Value = 0.2982772948921854
$ws.Range("P2").Value = 0.2359735829156887
$ws.Range("Q2").Value = 10.4939069243785
$ws.Range("R2").Value = 62.96344154627099
$ws.Range("S2").Value = 0.2982772948921854
$ws.Range("T2").Value = 0.2359735829156887

$ws.Range("N3").Value = 73.46982600000001
$ws.Range("O3").Value = 0.1616709822417395
$ws.Range("P3").Value = 0.1918520865636367
$ws.Range("Q3").Value = 5.687862499326
$ws.Range("R3").Value = 51.190762493934
$ws.Range("S3").Value = 0.1616709822417395
$ws.Range("T3").Value = 0.1918520865636367

$ws.Range("M4").Value = 20.755341
$ws.Range("N4").Value = 62.26602299999999
$ws.Range("O4").Value = 0.1370169176485697
$ws.Range("P4").Value = 0.1625955454769879
$ws.Range("Q4").Value = 4.820490213272999
$ws.Range("R4").Value = 43.38441191945699
$ws.Range("S4").Value = 0.1370169176485697
$ws.Range("T4").Value = 0.1625955454769879

$ws.Range("M5").Value = 26.3069545
$ws.Range("N5").Value = 52.613909
$ws.Range("O5").Value = 0.173666037012409
$ws.Range("P5").Value = 0.1373909368441856
$ws.Range("Q5").Value = 6.109869103488499
$ws.Range("R5").Value = 36.659214620931
$ws.Range("S5").Value = 0.173666037012409
$ws.Range("T5").Value = 0.1373909368441856

$ws.Range("M6").Value = 17.34473466666667
$ws.Range("N6").Value = 52.034204
$ws.Range("O6").Value = 0.1145017121838161
$ws.Range("P6").Value = 0.1358771505744131
$ws.Range("Q6").Value = 4.028366660537332
$ws.Range("R6").Value = 36.255299944836
$ws.Range("S6").Value = 0.1145017121838161
$ws.Range("T6").Value = 0.1358771505744131

$ws.Range("M7").Value = 17.400077
$ws.Range("N7").Value = 52.200231
$ws.Range("O7").Value = 0.1148670560212801
$ws.Range("P7").Value = 0.136310697625088
$ws.Range("Q7").Value = 4.041220083480999
$ws.Range("R7").Value = 36.37098075132899
$ws.Range("S7").Value = 0.1148670560212801
$ws.Range("T7").Value = 0.136310697625088
